# Dispatcher added to Initialization state and init application done
# successfully with orchestrator credential.
#
# Updates the Settings sheet: renames the orchestrator queue value and
# inserts new rows describing the UIDemo application / transaction data
# settings used by the Dispatcher/Performer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Queue name used by the process (Orchestrator queue).
$ws.Range("B2").Value = "UIDemoQueue"

# Insert 4 new configuration rows right above the existing
# "OrchestratorAssetFolder" row (old row 3), pushing everything below
# down by 4 rows.
$ws.Rows("3:6").Insert()
$ws.Rows("3:6").RowHeight = 14.25

$ws.Range("A3").Value = "UIDemoTransactionData_Path"
$ws.Range("B3").Value = "Data\Transactions.xlsx"

$ws.Range("A4").Value = "UIDemoApplication_Path"
$ws.Range("B4").Value = "UIDemo\UIDemo.exe"
$ws.Range("C4").Value = "Application path"

$ws.Range("A5").Value = "Input_SheetName"
$ws.Range("B5").Value = "Input"

# Row 6 is left blank (spacer row) with the same formatting as the rows
# above it, which Excel already applies automatically on insert.

# Move the active selection to match where the editor left off.
[void]$ws.Range("B14").Select()
